# Loan RBI, Variable Instalments
#
# The "Repayment schedule" sheet gets a new (blank) column inserted just
# before the existing "Late" column (column N), pushing Late / heading /
# Outstanding one column to the right (N->O, O->P, P->Q). Excel's normal
# "Insert Column" behaviour carries the formatting of the column to the
# left (M, "In Advance") into the freshly inserted column, so we copy that
# column's width across explicitly.
#
# The workbook's active sheet also moves from "Transactions" (last sheet)
# to "Repayment schedule", with the selection left on cell R8.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Remember the width of the column immediately to the left (M) so the new
# column inherits it, matching Excel's default insert-column behaviour.
$leftWidth = $ws.Columns("M:M").ColumnWidth

# Insert a new blank column at N, shifting the old N/O/P columns right.
$ws.Columns("N:N").Insert()
$ws.Columns("N:N").ColumnWidth = $leftWidth

# Make "Repayment schedule" the active sheet/tab, with R8 selected.
$ws.Activate()
[void]$ws.Range("R8").Select()
